$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C8").Value = 8915
$ws.Range("C9:C11").Value = 8678
$ws.Range("C12:C16").Value = 8241
$ws.Range("C17:C71").Value = 7670
$ws.Range("C72:C242").Value = 7586
$ws.Range("C243:C252").Value = 7312
